$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Event name as file name" - rename instructors to include surname "Słowianin"
$ws.Range("E10").Value = "Kajko Słowianin"
$ws.Range("E11").Value = "Kokosz Słowianin"

# Clear the duplicate/helper column J (keeps cell styles, drops values)
$ws.Range("J1").ClearContents()
$ws.Range("J2:J5").ClearContents()
$ws.Range("J6:J11").ClearContents()
$ws.Range("J12").ClearContents()

# Update the data validation formula on D1:D11 (case fix of broken name error)
$validation = $ws.Range("D1:D11").Validation
$validation.Formula1 = "#nazwa?"

# Move the active selection to L17
$ws.Range("L17").Select()
